$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = -20.44029999999998
$ws.Range("A7").Value = -21.94550000000001
$ws.Range("B7").Value = 4.741800000000003
$ws.Range("B15").Value = 4.913699999999997
$ws.Range("A16").Value = -21.66959999999999
$ws.Range("B21").Value = 10.53630000000001
$ws.Range("B22").Value = 10.2296
$ws.Range("B23").Value = 8.962700000000007
$ws.Range("A28").Value = -22.31879999999998
$ws.Range("A29").Value = -21.04989999999998
$ws.Range("A32").Value = -21.18569999999999
$ws.Range("B34").Value = 9.655400000000007
$ws.Range("A40").Value = -19.81549999999999
$ws.Range("B43").Value = 6.039800000000003
$ws.Range("B45").Value = 4.966100000000003
$ws.Range("B50").Value = 4.578699999999998
$ws.Range("B51").Value = 5.793199999999998
$ws.Range("A52").Value = -21.99099999999999
$ws.Range("A57").Value = -22.44940000000001
$ws.Range("A66").Value = -21.4489
$ws.Range("B66").Value = 5.3448
$ws.Range("B67").Value = 5.329400000000001
$ws.Range("B79").Value = 9.922800000000006
$ws.Range("B84").Value = 5.61
$ws.Range("B92").Value = 4.668299999999997
$ws.Range("B97").Value = 5.950399999999997
$ws.Range("A100").Value = -21.78949999999999

$wb.Save()
